$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -0.03550443442769693
$ws.Range("C2").Value = 0.4494482028570796
$ws.Range("D2").Value = 0.266414150275854
$ws.Range("E2").Value = 0.5161532236418309
$ws.Range("F2").Value = 0.5330038716810166
$ws.Range("G2").Value = 15

# Row 3
$ws.Range("B3").Value = 0.02310156321232606
$ws.Range("C3").Value = 0.3926748527515752
$ws.Range("D3").Value = 0.2150774634025442
$ws.Range("E3").Value = 0.463764448187379
$ws.Range("F3").Value = 0.4806736408029301
$ws.Range("G3").Value = 14

# Row 4
$ws.Range("B4").Value = -0.01768369464496384
$ws.Range("C4").Value = 0.5371479443576889
$ws.Range("D4").Value = 0.4289931411788038
$ws.Range("E4").Value = 0.6549756798376591
$ws.Range("F4").Value = 0.6814717874811226
$ws.Range("G4").Value = 13

# Row 5
$ws.Range("B5").Value = 0.04484032762995968
$ws.Range("C5").Value = 0.3900917434466585
$ws.Range("D5").Value = 0.2911825458959873
$ws.Range("E5").Value = 0.5396133299835979
$ws.Range("F5").Value = 0.5616584768642613
$ws.Range("G5").Value = 12

# Row 6
$ws.Range("B6").Value = 0.04618976013043215
$ws.Range("C6").Value = 0.3482013799665239
$ws.Range("D6").Value = 0.2021242622245639
$ws.Range("E6").Value = 0.4495823197419622
$ws.Range("F6").Value = 0.4690307507104655
$ws.Range("G6").Value = 11

# Row 7
$ws.Range("B7").Value = 0.01802609018820156
$ws.Range("C7").Value = 0.396353456599608
$ws.Range("D7").Value = 0.271167765387455
$ws.Range("E7").Value = 0.5207377126610431
$ws.Range("F7").Value = 0.5485767701364259
$ws.Range("G7").Value = 10

# Row 8
$ws.Range("B8").Value = -0.08230540368573748
$ws.Range("C8").Value = 0.4030200371981603
$ws.Range("D8").Value = 0.2441143047403328
$ws.Range("E8").Value = 0.4940792494532966
$ws.Range("F8").Value = 0.5167278209294698
$ws.Range("G8").Value = 9

# Row 9
$ws.Range("B9").Value = -0.03533701630595494
$ws.Range("C9").Value = 0.5782167566939017
$ws.Range("D9").Value = 0.5048288183583536
$ws.Range("E9").Value = 0.7105130669863529
$ws.Range("F9").Value = 0.7586304303617115
$ws.Range("G9").Value = 8

# Row 10
$ws.Range("B10").Value = 0.1068317712072389
$ws.Range("C10").Value = 0.4608479654118836
$ws.Range("D10").Value = 0.2863969179506793
$ws.Range("E10").Value = 0.5351606468628642
$ws.Range("F10").Value = 0.5664049249844474
$ws.Range("G10").Value = 7

# Row 11
$ws.Range("B11").Value = -0.005335926578848192
$ws.Range("C11").Value = 0.3488635011072845
$ws.Range("D11").Value = 0.1472458661147809
$ws.Range("E11").Value = 0.3837262906223405
$ws.Range("F11").Value = 0.4203104481247061
$ws.Range("G11").Value = 6
